$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.012.29"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.517.51"
$ws.Range("E3").Value = "  +3.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.57"
$ws.Range("E5").Value = "  +5.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.33"
$ws.Range("E6").Value = "  +4.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.566"
$ws.Range("E8").Value = "  +3.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.516.70"
$ws.Range("E9").Value = "  +2.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0995"
$ws.Range("E10").Value = "  +4.36%  "
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.333"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.961.59"
$ws.Range("E14").Value = "  +3.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.967.12"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.39"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("E17").Value = "  +3.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.517.55"
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.66"
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.25"
$ws.Range("E20").Value = "  +3.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.04"
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("E22").Value = "  +9.12%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.66"
$ws.Range("E24").Value = "  +3.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.411"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.51"
$ws.Range("E28").Value = "  +3.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0763"
$ws.Range("E29").Value = "  +5.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.46"
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.74"
$ws.Range("E31").Value = "  +5.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.20"
$ws.Range("E32").Value = "  +4.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.30"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.17"
$ws.Range("E36").Value = "  +2.43%  "
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.96"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.52"
$ws.Range("E39").Value = "  +4.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.57"
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.797"
$ws.Range("E41").Value = "  +4.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.49"
$ws.Range("E42").Value = "  +3.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "277.22"
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.06"
$ws.Range("E44").Value = "  +2.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "131.61"
$ws.Range("E45").Value = "  +9.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.593"
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0934"
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0511"
$ws.Range("E48").Value = "  +5.41%  "
$ws.Range("E49").Value = "  +5.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.06"
$ws.Range("E50").Value = "  +2.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.756.12"
$ws.Range("E51").Value = "  +3.11%  "
